# contratoh.xlsx update:
#  - row 27 (+556298529715 / 62 / 2024-07-09) is removed
#  - row 26 is overwritten with what used to be row 27's data
#  - sheet dimension shrinks from A1:C27 to A1:C26
#
# Row 26's text values (phone numbers, DDD codes, ISO dates) look numeric to
# Excel, so a plain `Range.Value = "..."` assignment would silently reinterpret
# them (dropping the leading "+", turning the date into a serial number, and
# bumping the cell's style to a new auto Text/Date format). To keep the cells
# as genuine literal text - matching the original inline-string cells and
# their existing style - we build each value as a text-literal formula
# (="...") and then collapse it down to a plain value in place via
# Copy + PasteSpecial(values only), which leaves the cell's existing
# formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last row (27); row 26 stays put, nothing shifts up automatically.
$ws.Rows.Item(27).Delete()

# Overwrite row 26 with the values that used to live in row 27.
$ws.Range("A26").Formula = "=""+556298529715"""
$ws.Range("B26").Formula = "=""62"""
$ws.Range("C26").Formula = "=""2024-07-09"""

# Convert those formulas to plain literal text values, preserving style/format.
$ws.Range("A26:C26").Copy()
$ws.Range("A26:C26").PasteSpecial(-4163)
